$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Corrected")
$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "adductName"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4131
Write-Host "done"
